$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the group N1 summary row: hackerrank score is now full (400/400)
# and the "Diem phat bieu" (speaking score) count increases from 6 to 10.
$ws.Range("B2").Value = "400/400"
$ws.Range("D2").Value = 10

# Update the selected cell shown when the workbook is reopened.
$ws.Range("F5").Select()
